$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 961.5333000000001
$ws.Range("I55").Value = 1625.8572
$ws.Range("J55").Value = 380.25
$ws.Range("K55").Value = 1625.8572
$ws.Range("L55").Value = 380.25
$ws.Range("M55").Value = -1411.8572
$ws.Range("N55").Value = -808.25

$ws.Range("H62").Value = 2810.1667
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2810.1667
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2810.1667
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4058.1667

$ws.Range("H65").Value = 2810.1667
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2810.1667
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 14050.8335
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -20290.8335

$ws.Range("H106").Value = 1986.9584
$ws.Range("I106").Value = 1290.5834
$ws.Range("K106").Value = 1290.5834
$ws.Range("M106").Value = -659.5834

$ws.Range("H107").Value = 755.1951
$ws.Range("I107").Value = 776.1613
$ws.Range("J107").Value = 690.2
$ws.Range("K107").Value = 776.1613
$ws.Range("L107").Value = 690.2
$ws.Range("M107").Value = 1143.8387
$ws.Range("N107").Value = -4530.2

$ws.Range("H113").Value = 4440.9287
$ws.Range("I113").Value = 3472.4546
$ws.Range("J113").Value = 7992
$ws.Range("K113").Value = 3472.4546
$ws.Range("L113").Value = 7992
$ws.Range("M113").Value = -218.4546
$ws.Range("N113").Value = -14500

$ws.Range("H136").Value = 26750
$ws.Range("J136").Value = 26750
$ws.Range("L136").Value = 26750
$ws.Range("N136").Value = -36950

$ws.Range("H137").Value = 2177542
$ws.Range("I137").Value = 3337356
$ws.Range("J137").Value = 2890.8125
$ws.Range("K137").Value = 10012068
$ws.Range("L137").Value = 8672.4375
$ws.Range("M137").Value = -10009518
$ws.Range("N137").Value = -13772.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4808889
$ws.Range("I2").Value = 8334187
$ws.Range("K2").Value = 8334187
$ws.Range("M2").Value = -8334074

$ws.Range("H32").Value = 2950.42
$ws.Range("I32").Value = 2950.42
$ws.Range("K32").Value = 2950.42
$ws.Range("M32").Value = -2663.42

$ws.Range("H74").Value = 1703.4642
$ws.Range("I74").Value = 1199.7
$ws.Range("J74").Value = 2962.875
$ws.Range("K74").Value = 1199.7
$ws.Range("L74").Value = 2962.875
$ws.Range("M74").Value = -325.7
$ws.Range("N74").Value = -4710.875

$ws.Range("H77").Value = 1703.4642
$ws.Range("I77").Value = 1199.7
$ws.Range("J77").Value = 2962.875
$ws.Range("K77").Value = 5998.5
$ws.Range("L77").Value = 14814.375
$ws.Range("M77").Value = -1630.5
$ws.Range("N77").Value = -23550.375

$ws.Range("H110").Value = 2094
$ws.Range("I110").Value = 620.7143
$ws.Range("J110").Value = 4385.778
$ws.Range("K110").Value = 620.7143
$ws.Range("L110").Value = 4385.778
$ws.Range("M110").Value = 1424.2857
$ws.Range("N110").Value = -8475.778

$ws.Range("H116").Value = 4808889
$ws.Range("I116").Value = 8334187
$ws.Range("K116").Value = 8334187
$ws.Range("M116").Value = -8331893

$ws.Range("H132").Value = 2051.5166
$ws.Range("I132").Value = 1691.3
$ws.Range("K132").Value = 5073.9
$ws.Range("M132").Value = -2543.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4808889
$ws.Range("I3").Value = 8334187
$ws.Range("K3").Value = 8334187
$ws.Range("M3").Value = -8334073

$ws.Range("H99").Value = 3570.4211
$ws.Range("I99").Value = 3231.2856
$ws.Range("K99").Value = 3231.2856
$ws.Range("M99").Value = -1733.2856

$ws.Range("H105").Value = 1507.5834
$ws.Range("I105").Value = 1388
$ws.Range("K105").Value = 1388
$ws.Range("M105").Value = 359

$ws.Range("H107").Value = 2175.3333
$ws.Range("I107").Value = 1099.6666
$ws.Range("K107").Value = 1099.6666
$ws.Range("M107").Value = 820.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 69174.836
$ws.Range("I23").Value = 7504.5
$ws.Range("J23").Value = 100010
$ws.Range("K23").Value = 7504.5
$ws.Range("L23").Value = 100010
$ws.Range("M23").Value = -7264.5
$ws.Range("N23").Value = -100490

$ws.Range("H27").Value = 69174.836
$ws.Range("I27").Value = 7504.5
$ws.Range("J27").Value = 100010
$ws.Range("K27").Value = 7504.5
$ws.Range("L27").Value = 100010
$ws.Range("M27").Value = -7312.5
$ws.Range("N27").Value = -100394

$ws.Range("H134").Value = 15627861
$ws.Range("I134").Value = 22729530
$ws.Range("K134").Value = 68188590
$ws.Range("M134").Value = -68186055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1133.3334
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1133.3334
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3400.0002
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5772.0002

$ws.Range("H89").Value = 1133.3334
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1133.3334
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 10200.0006
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -22056.0006

$ws.Range("H97").Value = 1709
$ws.Range("J97").Value = 2187.375
$ws.Range("L97").Value = 6562.125
$ws.Range("N97").Value = -7554.125

$ws.Range("H122").Value = 924.2
$ws.Range("I122").Value = 461.2
$ws.Range("J122").Value = 1232.8667
$ws.Range("K122").Value = 4150.8
$ws.Range("L122").Value = 11095.8003
$ws.Range("M122").Value = -1700.8
$ws.Range("N122").Value = -15995.8003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 9000
$ws.Range("J27").Value = 9000
$ws.Range("L27").Value = 9000
$ws.Range("N27").Value = -9332

$ws.Range("H46").Value = 12147.6
$ws.Range("J46").Value = 12934.5
$ws.Range("L46").Value = 12934.5
$ws.Range("N46").Value = -13246.5

$ws.Range("H102").Value = 54343.1
$ws.Range("I102").Value = 3395.9
$ws.Range("J102").Value = 105290.3
$ws.Range("K102").Value = 3395.9
$ws.Range("L102").Value = 105290.3
$ws.Range("M102").Value = -1773.9
$ws.Range("N102").Value = -108534.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1890
$ws.Range("I46").Value = 518.5714
$ws.Range("J46").Value = 2850
$ws.Range("K46").Value = 518.5714
$ws.Range("L46").Value = 2850
$ws.Range("M46").Value = -330.5714
$ws.Range("N46").Value = -3226

$ws.Range("H68").Value = 2523.3157
$ws.Range("I68").Value = 1174.2858
$ws.Range("J68").Value = 6300.6
$ws.Range("K68").Value = 1174.2858
$ws.Range("L68").Value = 6300.6
$ws.Range("M68").Value = -425.2858000000001
$ws.Range("N68").Value = -7798.6

$ws.Range("H71").Value = 2523.3157
$ws.Range("I71").Value = 1174.2858
$ws.Range("J71").Value = 6300.6
$ws.Range("K71").Value = 5871.429
$ws.Range("L71").Value = 31503
$ws.Range("M71").Value = -2127.429
$ws.Range("N71").Value = -38991

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H50").Value = 5000
$ws.Range("J50").Value = 5000
$ws.Range("L50").Value = 5000
$ws.Range("N50").Value = -6262

$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 5000
$ws.Range("K52").Value = 5000
$ws.Range("M52").Value = -4774

$ws.Range("H53").Value = 20084
$ws.Range("J53").Value = 20084
$ws.Range("L53").Value = 20084
$ws.Range("N53").Value = -21298

$ws.Range("H54").Value = 10593.9
$ws.Range("I54").Value = 8000
$ws.Range("J54").Value = 11242.375
$ws.Range("K54").Value = 8000
$ws.Range("L54").Value = 11242.375
$ws.Range("M54").Value = -7480
$ws.Range("N54").Value = -12282.375

$ws.Range("H61").Value = 23333.334
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H96").Value = 1136.6364
$ws.Range("I96").Value = 1250.5
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1250.5
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 122.5
$ws.Range("N96").Value = -3746

$ws.Range("H97").Value = 27722.223
$ws.Range("J97").Value = 28312.5
$ws.Range("L97").Value = 28312.5
$ws.Range("N97").Value = -30294.5

$ws.Range("H101").Value = 19575.125
$ws.Range("J101").Value = 19575.125
$ws.Range("L101").Value = 19575.125
$ws.Range("N101").Value = -26065.125

$ws.Range("H104").Value = 34439
$ws.Range("J104").Value = 34439
$ws.Range("L104").Value = 34439
$ws.Range("N104").Value = -41427
